# towards new measurement integration
# Fill in newly-recorded GPS coordinates for a few samples (rows 12-14, 19)
# and normalize number formatting on the coordinate columns for rows 15-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New coordinate readings for MSTW-21-11 (row 12)
$ws.Range("F12").Value = "3'08.519"
$ws.Range("G12").Value = "35'51.721"

# New coordinate readings for MSTW-21-12 (row 13)
$ws.Range("F13").Value = "3'08.520"
$ws.Range("G13").Value = "35'51.722"

# New coordinate readings for MSTW-21-13 (row 14)
$ws.Range("F14").Value = "3'08.521"
$ws.Range("G14").Value = "35'51.723"

# New coordinate readings for MSTW-21-18 (row 19)
$ws.Range("F19").Value = "3'08.521"
$ws.Range("G19").Value = "35'51.723"

# Normalize the numeric lat/long columns (rows 15-18) to the same 5-decimal
# number format already used elsewhere in the sheet.
$ws.Range("G15").NumberFormat = "0.00000"
$ws.Range("F16:G16").NumberFormat = "0.00000"
$ws.Range("F17:G17").NumberFormat = "0.00000"
$ws.Range("F18:G18").NumberFormat = "0.00000"

# Move the active selection to reflect where editing left off.
$ws.Range("I16").Select() | Out-Null
